$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values for rows 2-5 across columns A:AH ---
$arr = New-Object 'object[,]' 4,34
$arr[0,0] = 45169.50694444445
$arr[0,1] = 20.178
$arr[0,2] = 13.652
$arr[0,3] = 4.068
$arr[0,4] = 42.752
$arr[0,5] = 34.691
$arr[0,6] = 15.879
$arr[0,7] = 51.007
$arr[0,8] = 24.432
$arr[0,9] = 10.251
$arr[0,10] = 15.607
$arr[0,11] = 16.868
$arr[0,12] = 17.597
$arr[0,13] = 5.069
$arr[0,14] = 15.79
$arr[0,15] = 22.094
$arr[0,16] = 13.41
$arr[0,17] = 3.498
$arr[0,18] = 2.451
$arr[0,19] = 232.996
$arr[0,20] = 43.923
$arr[0,21] = 14.575
$arr[0,22] = 29.082
$arr[0,23] = 14.996
$arr[0,24] = 3.13
$arr[0,25] = 25.159
$arr[0,26] = 12.874
$arr[0,27] = 11.654
$arr[0,28] = 13.651
$arr[0,29] = 17.331
$arr[0,30] = 3.457
$arr[0,31] = 45.214
$arr[0,32] = 8.105
$arr[0,33] = 18.222
$arr[1,0] = 45169.51388888889
$arr[1,1] = 9.609
$arr[1,2] = 6.511
$arr[1,3] = 1.581
$arr[1,4] = 20.565
$arr[1,5] = 16.566
$arr[1,6] = 7.562
$arr[1,7] = 31.653
$arr[1,8] = 11.635
$arr[1,9] = 4.908
$arr[1,10] = 7.264
$arr[1,11] = 8.228999999999999
$arr[1,12] = 8.516
$arr[1,13] = 2.418
$arr[1,14] = 7.519
$arr[1,15] = 10.522
$arr[1,16] = 6.618
$arr[1,17] = 1.509
$arr[1,18] = 0.885
$arr[1,19] = 107.143
$arr[1,20] = 21.161
$arr[1,21] = 6.941
$arr[1,22] = 13.895
$arr[1,23] = 7.251
$arr[1,24] = 1.576
$arr[1,25] = 14.794
$arr[1,26] = 6.131
$arr[1,27] = 5.655
$arr[1,28] = 6.607
$arr[1,29] = 8.411
$arr[1,30] = 1.266
$arr[1,31] = 28.943
$arr[1,32] = 3.8
$arr[1,33] = 8.678000000000001
$arr[2,0] = 45169.52083333334
$arr[2,1] = 14.413
$arr[2,2] = 10.367
$arr[2,3] = 1.254
$arr[2,4] = 31.178
$arr[2,5] = 25.429
$arr[2,6] = 11.342
$arr[2,7] = 42.99
$arr[2,8] = 17.452
$arr[2,9] = 7.616
$arr[2,10] = 11.291
$arr[2,11] = 12.527
$arr[2,12] = 13.109
$arr[2,13] = 3.622
$arr[2,14] = 11.279
$arr[2,15] = 15.94
$arr[2,16] = 9.657
$arr[2,17] = 1.04
$arr[2,18] = 0.765
$arr[2,19] = 164.344
$arr[2,20] = 31.506
$arr[2,21] = 10.411
$arr[2,22] = 21.011
$arr[2,23] = 11.07
$arr[2,24] = 1.887
$arr[2,25] = 20.715
$arr[2,26] = 9.196
$arr[2,27] = 8.263999999999999
$arr[2,28] = 9.694000000000001
$arr[2,29] = 13.015
$arr[2,30] = 0.773
$arr[2,31] = 38.87
$arr[2,32] = 5.81
$arr[2,33] = 13.016
$arr[3,0] = 45169.52777777778
$arr[3,1] = 19.22
$arr[3,2] = 14.1
$arr[3,3] = 1.2
$arr[3,4] = 41.69
$arr[3,5] = 34.18
$arr[3,6] = 15.12
$arr[3,7] = 57.81
$arr[3,8] = 23.27
$arr[3,9] = 10.29
$arr[3,10] = 15.27
$arr[3,11] = 16.75
$arr[3,12] = 17.61
$arr[3,13] = 4.83
$arr[3,14] = 15.04
$arr[3,15] = 21.37
$arr[3,16] = 12.71
$arr[3,17] = 0.84
$arr[3,18] = 0.8
$arr[3,19] = 221.58
$arr[3,20] = 42.04
$arr[3,21] = 13.88
$arr[3,22] = 28.22
$arr[3,23] = 14.85
$arr[3,24] = 2.28
$arr[3,25] = 28.09
$arr[3,26] = 12.26
$arr[3,27] = 10.91
$arr[3,28] = 12.81
$arr[3,29] = 17.51
$arr[3,30] = 0.5600000000000001
$arr[3,31] = 52.43
$arr[3,32] = 7.81
$arr[3,33] = 17.35
$ws.Range("A2:AH5").Value = $arr

# --- Delete row 6 (data now only spans rows 1-5) ---
$ws.Rows.Item(6).Delete()

# --- Widen columns C, J, X, AA, AB from 7 to 8 characters ---
$ws.Columns.Item(3).ColumnWidth = 7.17
$ws.Columns.Item(10).ColumnWidth = 7.17
$ws.Columns.Item(24).ColumnWidth = 7.17
$ws.Columns.Item(27).ColumnWidth = 7.17
$ws.Columns.Item(28).ColumnWidth = 7.17

"Done"